$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed values
# Row 2
$ws.Range("D2").Value = 214474
$ws.Range("E2").Value = 19591
$ws.Range("F2").Value = 19591
$ws.Range("G2").Value = 19014
$ws.Range("H2").Value = 14151
$ws.Range("I2").Value = 14007
$ws.Range("J2").Value = 144
$ws.Range("K2").Value = 3083557
$ws.Range("L2").Value = 2808430
$ws.Range("M2").Value = 275127
$ws.Range("N2").Value = 273151
$ws.Range("O2").Value = 1976
$ws.Range("P2").Value = 19318
$ws.Range("Q2").Value = 32767
$ws.Range("R2").Value = -16843
$ws.Range("S2").Value = -7551
$ws.Range("T2").Value = 2020
$ws.Range("V2").Value = 411094
$ws.Range("W2").Value = 9.130000000000001
$ws.Range("X2").Value = 6.6
$ws.Range("Y2").Value = 5.26
$ws.Range("Z2").Value = 0.47
$ws.Range("AA2").Value = 1020.78
$ws.Range("AB2").Value = 1324.23
$ws.Range("AC2").Value = 3626
$ws.Range("AD2").Value = 9.970000000000001
$ws.Range("AE2").Value = 70700
$ws.Range("AF2").Value = 0.51
$ws.Range("AG2").Value = 780
$ws.Range("AH2").Value = 2.16
$ws.Range("AI2").Value = 21.51
$ws.Range("AJ2").Value = 386351693

# Row 3
$ws.Range("D3").Value = 222333
$ws.Range("E3").Value = 18211
$ws.Range("F3").Value = 18211
$ws.Range("G3").Value = 21647
$ws.Range("H3").Value = 17273
$ws.Range("I3").Value = 16983
$ws.Range("J3").Value = 290
$ws.Range("K3").Value = 3290655
$ws.Range("L3").Value = 3001627
$ws.Range("M3").Value = 289027
$ws.Range("N3").Value = 286806
$ws.Range("O3").Value = 2221
$ws.Range("P3").Value = 19318
$ws.Range("Q3").Value = 21936
$ws.Range("R3").Value = -50753
$ws.Range("S3").Value = 32553
$ws.Range("T3").Value = 2292
$ws.Range("V3").Value = 449048
$ws.Range("W3").Value = 8.19
$ws.Range("X3").Value = 7.77
$ws.Range("Y3").Value = 6.07
$ws.Range("Z3").Value = 0.54
$ws.Range("AA3").Value = 1038.53
$ws.Range("AB3").Value = 1396.19
$ws.Range("AC3").Value = 4396
$ws.Range("AD3").Value = 7.54
$ws.Range("AE3").Value = 74234
$ws.Range("AF3").Value = 0.45
$ws.Range("AG3").Value = 980
$ws.Range("AH3").Value = 2.96
$ws.Range("AI3").Value = 22.29
$ws.Range("AJ3").Value = 386351693

# Row 4
$ws.Range("D4").Value = 253558
$ws.Range("E4").Value = 16769
$ws.Range("F4").Value = 16769
$ws.Range("G4").Value = 26287
$ws.Range("H4").Value = 21902
$ws.Range("I4").Value = 21437
$ws.Range("J4").Value = 464
$ws.Range("K4").Value = 3756737
$ws.Range("L4").Value = 3444123
$ws.Range("M4").Value = 312614
$ws.Range("N4").Value = 309980
$ws.Range("O4").Value = 2634
$ws.Range("P4").Value = 20906
$ws.Range("Q4").Value = 11250
$ws.Range("R4").Value = -44385
$ws.Range("S4").Value = 31813
$ws.Range("T4").Value = 3972
$ws.Range("V4").Value = 494778
$ws.Range("W4").Value = 6.61
$ws.Range("X4").Value = 8.640000000000001
$ws.Range("Y4").Value = 7.18
$ws.Range("Z4").Value = 0.62
$ws.Range("AA4").Value = 1101.72
$ws.Range("AB4").Value = 1429.9
$ws.Range("AC4").Value = 5459
$ws.Range("AD4").Value = 7.84
$ws.Range("AE4").Value = 77815
$ws.Range("AF4").Value = 0.55
$ws.Range("AG4").Value = 1250
$ws.Range("AH4").Value = 2.92
$ws.Range("AI4").Value = 23.23
$ws.Range("AJ4").Value = 418111537

# Row 5
$ws.Range("D5").Value = 392293
$ws.Range("E5").Value = 40153
$ws.Range("F5").Value = 40153
$ws.Range("G5").Value = 41384
$ws.Range("H5").Value = 33435
$ws.Range("I5").Value = 33114
$ws.Range("J5").Value = 320
$ws.Range("K5").Value = 4367856
$ws.Range("L5").Value = 4027408
$ws.Range("M5").Value = 340448
$ws.Range("N5").Value = 340387
$ws.Range("O5").Value = 61
$ws.Range("P5").Value = 20906
$ws.Range("Q5").Value = -48360
$ws.Range("R5").Value = -87294
$ws.Range("S5").Value = 146887
$ws.Range("T5").Value = 2984
$ws.Range("V5").Value = 618388
$ws.Range("W5").Value = 10.23
$ws.Range("X5").Value = 8.52
$ws.Range("Y5").Value = 10.18
$ws.Range("Z5").Value = 0.82
$ws.Range("AA5").Value = 1182.97
$ws.Range("AB5").Value = 1564.67
$ws.Range("AC5").Value = 7920
$ws.Range("AD5").Value = 8.01
$ws.Range("AE5").Value = 85302
$ws.Range("AF5").Value = 0.74
$ws.Range("AG5").Value = 1920
$ws.Range("AH5").Value = 3.03
$ws.Range("AI5").Value = 23.15
$ws.Range("AJ5").Value = 418111537

# Row 6
$ws.Range("D6").Value = 420271
$ws.Range("E6").Value = 42675
$ws.Range("F6").Value = 42675
$ws.Range("G6").Value = 43015
$ws.Range("H6").Value = 30619
$ws.Range("I6").Value = 30612
$ws.Range("K6").Value = 4795883
$ws.Range("L6").Value = 4438753
$ws.Range("M6").Value = 357130
$ws.Range("N6").Value = 357039
$ws.Range("P6").Value = 20906
$ws.Range("Q6").Value = -89576
$ws.Range("R6").Value = -44803
$ws.Range("S6").Value = 117438
$ws.Range("T6").Value = 4523
$ws.Range("V6").Value = 732480
$ws.Range("W6").Value = 10.15
$ws.Range("X6").Value = 7.29
$ws.Range("Y6").Value = 8.779999999999999
$ws.Range("Z6").Value = 0.67
$ws.Range("AA6").Value = 1242.89
$ws.Range("AB6").Value = 1654.63
$ws.Range("AC6").Value = 7321
$ws.Range("AD6").Value = 6.35
$ws.Range("AE6").Value = 90264
$ws.Range("AF6").Value = 0.52
$ws.Range("AG6").Value = 1920
$ws.Range("AH6").Value = 4.13
$ws.Range("AI6").Value = 24.82
$ws.Range("AJ6").Value = 418111537

# Row 7
$ws.Range("E7").Value = 45640
$ws.Range("G7").Value = 45499
$ws.Range("H7").Value = 33403
$ws.Range("I7").Value = 33347
$ws.Range("K7").Value = 5104758
$ws.Range("L7").Value = 4716694
$ws.Range("M7").Value = 388063
$ws.Range("N7").Value = 380485
$ws.Range("P7").Value = 20900
$ws.Range("Y7").Value = 9.039999999999999
$ws.Range("Z7").Value = 0.68
$ws.Range("AA7").Value = 1215.44
$ws.Range("AC7").Value = 7978
$ws.Range("AD7").Value = 5.51
$ws.Range("AE7").Value = 97652
$ws.Range("AF7").Value = 0.45
$ws.Range("AG7").Value = 2144
$ws.Range("AH7").Value = 4.88
$ws.Range("AI7").Value = 26.73

# Row 8
$ws.Range("E8").Value = 45847
$ws.Range("G8").Value = 45406
$ws.Range("H8").Value = 33435
$ws.Range("I8").Value = 33082
$ws.Range("K8").Value = 5298722
$ws.Range("L8").Value = 4885786
$ws.Range("M8").Value = 412995
$ws.Range("N8").Value = 406910
$ws.Range("P8").Value = 20900
$ws.Range("Y8").Value = 8.4
$ws.Range("Z8").Value = 0.64
$ws.Range("AA8").Value = 1183.01
$ws.Range("AC8").Value = 7956
$ws.Range("AD8").Value = 5.52
$ws.Range("AE8").Value = 104434
$ws.Range("AF8").Value = 0.42
$ws.Range("AG8").Value = 2237
$ws.Range("AH8").Value = 5.09
$ws.Range("AI8").Value = 28.12

# Row 9
$ws.Range("E9").Value = 48130
$ws.Range("G9").Value = 47092
$ws.Range("H9").Value = 34695
$ws.Range("I9").Value = 33468
$ws.Range("K9").Value = 5514696
$ws.Range("L9").Value = 5074266
$ws.Range("M9").Value = 440430
$ws.Range("N9").Value = 435111
$ws.Range("P9").Value = 20898
$ws.Range("Y9").Value = 7.95
$ws.Range("Z9").Value = 0.64
$ws.Range("AA9").Value = 1152.12
$ws.Range("AC9").Value = 8049
$ws.Range("AD9").Value = 5.46
$ws.Range("AE9").Value = 111672
$ws.Range("AF9").Value = 0.39
$ws.Range("AG9").Value = 2373
$ws.Range("AH9").Value = 5.4
$ws.Range("AI9").Value = 29.48

# Clear cells that were removed entirely in the update
$clearAddrs = @("U2","U3","U4","U5","U6","D7","Q7","R7","S7","T7","U7","W7","X7","D8","Q8","R8","S8","T8","U8","W8","X8","D9","Q9","R9","S9","T9","U9","W9","X9")
foreach ($addr in $clearAddrs) {
    $ws.Range($addr).ClearContents()
}
